$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Piña" (Terminal
# Hortofrutícola Agro Chillán). It belongs right after the current row 336,
# so insert a fresh row at 337 - this pushes the old rows 337:364 down to
# 338:365 (and grows the sheet's used range to A1:T365), then fill the new
# row with its own data.
$ws.Rows.Item(337).Insert()

$ws.Range("A337").Value = 7
$ws.Range("B337").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C337").Value = "Ñuble"
$ws.Range("D337").Value = 45166
$ws.Range("E337").Value = 16
$ws.Range("F337").Value = "Fruta"
$ws.Range("G337").Value = 100108
$ws.Range("H337").Value = "Tropicales y subtropicales"
$ws.Range("I337").Value = 100108005
$ws.Range("J337").Value = "Piña"
$ws.Range("K337").Value = "Caramelo"
$ws.Range("L337").Value = "Segunda"
$ws.Range("M337").Value = 60
$ws.Range("N337").Value = 22000
$ws.Range("O337").Value = 22000
$ws.Range("P337").Value = 22000
$ws.Range("Q337").Value = "$/caja 14 unidades"
$ws.Range("R337").Value = "Ecuador"
$ws.Range("S337").Value = 1571
$ws.Range("T337").Value = 14
